$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$c = $ws.Range("D22")
$c.ClearFormats()
$c.Interior.Color = 15518118
$c.WrapText = $true
$c.HorizontalAlignment = -4108
